$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new log row (row 21) to match the existing "Source Branch",
# "Author", "Action", "Comment", "Date", "Change ID" columns (A-F).
$ws.Range("A21").Value = "edit1"
$ws.Range("B21").Value = "riya-morankar"
$ws.Range("C21").Value = "Merged"
$ws.Range("D21").Value = "change1"

# Keep the Date column as literal text (matching the rest of the sheet)
# instead of letting Excel auto-convert the string into a date serial.
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2025-06-23"

$ws.Range("F21").Value = "1cf76dd27c2ae0ba57e39ad101ab17a10e2d488a"
